$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-08-16 Wednesday" "2023-08-17 Thursday"

Replace-Text "95×23=" "95×87="
Replace-Text "43×68=" "40×27="
Replace-Text "29×22=" "77×48="
Replace-Text "60×66=" "65×40="
Replace-Text "48×33=" "31×63="
Replace-Text "19×43=" "98×49="
Replace-Text "81×99=" "28×60="
Replace-Text "85×52=" "14×39="
Replace-Text "72×68=" "19×47="
Replace-Text "70×75=" "25×98="
Replace-Text "80×97=" "51×58="
Replace-Text "37×13=" "36×64="
Replace-Text "13×70=" "25×31="
Replace-Text "32×77=" "88×17="
Replace-Text "17×92=" "40×28="
Replace-Text "29×32=" "75×37="
Replace-Text "43×41=" "23×49="
Replace-Text "40×41=" "97×42="
Replace-Text "75×54=" "82×73="
Replace-Text "38×50=" "86×44="
Replace-Text "60×57=" "82×83="
Replace-Text "89×60=" "99×33="
Replace-Text "95×15=" "34×22="
Replace-Text "26×34=" "80×55="
Replace-Text "67×64=" "86×21="
